$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 243, shifting existing rows 243:347 down to 244:348
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with its data
$ws.Range("A243").Value = 8
$ws.Range("B243").Value = "Terminal La Palmera de La Serena"
$ws.Range("C243").Value = "Coquimbo"
$ws.Range("D243").Value = 44704
$ws.Range("E243").Value = 4
$ws.Range("F243").Value = 100114013
$ws.Range("G243").Value = "Zanahoria"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 600
$ws.Range("K243").Value = 6000
$ws.Range("L243").Value = 7000
$ws.Range("M243").Value = 6500
$ws.Range("N243").Value = "`$/saco 20 kilos"
$ws.Range("O243").Value = "Provincia del Elquí"
$ws.Range("P243").Value = 325
$ws.Range("Q243").Value = 20
$ws.Range("R243").Value = "Hortaliza"
